# Generate Report for Handback
# Renames the handed-back source file UUIDs (and their derived .xlf
# correspondence files + timestamps) across the Overview / zh-cn / de-de
# sheets, keeping each hyperlink's underlying Address untouched and only
# updating the visible display text (and therefore the cell's text value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New values (per the commit's regenerated handback report)
# ---------------------------------------------------------------------
$oldFile1 = "93bd694f-a1a8-4d8e-9745-eb97625df1f4.md"
$newFile1 = "f5b2ccdf-10b2-4339-94e9-4b45c024f529.md"

$oldFile2 = "a8eed73e-ebd3-4b54-8e2a-535cfbadb727.md"
$newFile2 = "ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57.md"

$newXlfZh = "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.zh-cn.xlf"
$newXlfDe = "f5b2ccdf-10b2-4339-94e9-4b45c024f529.ee3dcb8c376b760ac1fd1be7b860a5fba442b146.de-de.xlf"

$newHandoffZh = "2016-03-23 17:16:25"
$newHandbackZh = "2016-03-23 17:17:06"

$newHandoffDe = "2016-03-23 17:16:29"
$newHandbackDe = "2016-03-23 17:17:14"

# ---------------------------------------------------------------------
# Overview sheet: just the two file-name hyperlinks
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsOverview.Range("A3").Hyperlinks.Item(1).TextToDisplay = $newFile2

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsZh.Range("F2").Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsZh.Range("D2").Hyperlinks.Item(1).TextToDisplay = $newXlfZh
$wsZh.Range("G2").Hyperlinks.Item(1).TextToDisplay = $newXlfZh
$wsZh.Range("E2").Value = $newHandoffZh
$wsZh.Range("H2").Value = $newHandbackZh

$wsZh.Range("A3").Hyperlinks.Item(1).TextToDisplay = $newFile2
$wsZh.Range("F3").Hyperlinks.Item(1).TextToDisplay = $newFile2
$wsZh.Range("D3").Hyperlinks.Item(1).TextToDisplay = $newXlfZh
$wsZh.Range("G3").Hyperlinks.Item(1).TextToDisplay = $newXlfZh
$wsZh.Range("E3").Value = $newHandoffZh
$wsZh.Range("H3").Value = $newHandbackZh

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsDe.Range("F2").Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsDe.Range("D2").Hyperlinks.Item(1).TextToDisplay = $newXlfDe
$wsDe.Range("G2").Hyperlinks.Item(1).TextToDisplay = $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDe
$wsDe.Range("H2").Value = $newHandbackDe

$wsDe.Range("A3").Hyperlinks.Item(1).TextToDisplay = $newFile2
$wsDe.Range("F3").Hyperlinks.Item(1).TextToDisplay = $newFile2
$wsDe.Range("D3").Hyperlinks.Item(1).TextToDisplay = $newXlfDe
$wsDe.Range("G3").Hyperlinks.Item(1).TextToDisplay = $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDe
$wsDe.Range("H3").Value = $newHandbackDe
